$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing figures for row 8 (ABASTOS 11 SUR / ODELPA columns)
$ws.Range("B8").Value = 986
$ws.Range("C8").Value = 284
$ws.Range("L8").Value = 1392
$ws.Range("M8").Value = 287

# Move the view/selection: scroll so column B is the leftmost visible
# column and select M9 (next entry row, last column)
$ws.Range("M9").Select()
$excel.ActiveWindow.ScrollColumn = 2
